# Weekly price-sheet update: a new record (week of 2023-03-24) is
# prepended to the "Feria Lagunitas de Puerto Montt - Cebollín" price
# history. All existing records from row 301 down shift one row lower
# (old row 301 -> new row 302, ..., old row 426 -> new row 427), and the
# new data lands in the now-vacant row 301.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 301; this pushes rows
# 301..426 down to 302..427 (and grows the sheet dimension to A1:R427).
$ws.Rows.Item(301).Insert()

# Populate the newly inserted row 301 with the new weekly record.
$ws.Range("A301").Value = 4
$ws.Range("B301").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C301").Value = "Los Lagos"
$ws.Range("D301").Value = 45009
$ws.Range("E301").Value = 10
$ws.Range("F301").Value = 100112037
$ws.Range("G301").Value = "Cebollín"
$ws.Range("H301").Value = "Sin especificar"
$ws.Range("I301").Value = "Primera"
$ws.Range("J301").Value = 180
$ws.Range("K301").Value = 6500
$ws.Range("L301").Value = 7000
$ws.Range("M301").Value = 6750
$ws.Range("N301").Value = "`$/paquete 36 unidades"
$ws.Range("O301").Value = "Región Metropolitana"
$ws.Range("P301").Value = 188
$ws.Range("Q301").Value = 36
$ws.Range("R301").Value = "Hortaliza"
